# Add a missing "9:00 - 10:00" interval row to the schedule.
# This inserts a new row above the existing row 3 ("10:00 - 11:00"),
# shifting all subsequent rows down by one, and fills the new row's
# A-cell with the missing interval label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 3 - everything from row 3 down shifts to row+1.
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row with the missing time interval.
$ws.Range("A3").Value = "9:00 - 10:00"

# Move the active selection to the newly inserted cell, matching the
# post-edit workbook state.
$ws.Range("A3").Select()
